$d = $word.ActiveDocument

# The document currently ends with a single empty ListParagraph bullet at
# ilvl=1 (numId=15). We turn it into a "Python code" bullet at ilvl=0 and
# then append six more bullets underneath it describing how to implement
# kfolds cross validation and score comparison with cross_val_score.

$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 1   # w:ilvl 0
$p.Range.Text = "Python code"
$p.Range.NoProofing = 1

$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.ListFormat.ListLevelNumber = 2   # w:ilvl 1
$p2.Range.Text = "Implement kfolds"
$p2.Range.NoProofing = 1

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.ListFormat.ListLevelNumber = 3   # w:ilvl 2
$p3.Range.Text = "from sklearn.model_selection import StratifiedKFold"
$p3.Range.NoProofing = 1

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.ListFormat.ListLevelNumber = 3   # w:ilvl 2
$p4.Range.Text = "folds = StratifiedKFold(n_splits=3)"
$p4.Range.NoProofing = 1

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Range.ListFormat.ListLevelNumber = 2   # w:ilvl 1
$p5.Range.Text = "Automatically get a score that compares different ML models"
$p5.Range.NoProofing = 1

$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Last
$p6.Range.ItalicBi = 1
$p6.Range.ListFormat.ListLevelNumber = 3   # w:ilvl 2
$p6.Range.Text = "from sklearn.model_selection import cross_val_score"
$p6.Range.Italic = 1
$p6.Range.ItalicBi = 1
$p6.Range.NoProofing = 1

$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Last
$p7.Range.ItalicBi = 1
$p7.Range.ListFormat.ListLevelNumber = 3   # w:ilvl 2
$p7.Range.Text = "cross_val_score(RandomForestClassifier(), digits.data, digits.target)"
$p7.Range.Italic = 1
$p7.Range.ItalicBi = 1
$p7.Range.NoProofing = 1
